# Update the workbook with the new TPM-derived NATMI values.
# Rows 2-4 (Sending=ECs) get refreshed numbers and a corrected
# Sending/Ligand/Receptor/Target column mapping; rows 5-10 are new
# (Sending=FAPs and Sending=MuSCs blocks), completing the full 3x3
# Sending-cluster x Target-cluster grid for the Wnt1->Ror2 pair.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt1"
$ws.Range("C2").Value = "Ror2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.09370099999999999
$ws.Range("H2").Value = 0.281103
$ws.Range("I2").Value = 0.6634967391997054
$ws.Range("J2").Value = 0.6634967391997054
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.11549
$ws.Range("N2").Value = 0.34647
$ws.Range("O2").Value = 0.01449407350231777
$ws.Range("P2").Value = 0.01449407350231777
$ws.Range("Q2").Value = 0.01082152849
$ws.Range("R2").Value = 0.09739375640999999
$ws.Range("S2").Value = 0.009616770506508691
$ws.Range("T2").Value = 0.009616770506508693

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt1"
$ws.Range("C3").Value = "Ror2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.09370099999999999
$ws.Range("H3").Value = 0.281103
$ws.Range("I3").Value = 0.6634967391997054
$ws.Range("J3").Value = 0.6634967391997054
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 7.692787333333334
$ws.Range("N3").Value = 23.078362
$ws.Range("O3").Value = 0.9654500393716549
$ws.Range("P3").Value = 0.965450039371655
$ws.Range("Q3").Value = 0.7208218659206667
$ws.Range("R3").Value = 6.487396793286
$ws.Range("S3").Value = 0.6405729529833202
$ws.Range("T3").Value = 0.6405729529833203

# Row 4: ECs -> MuSCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt1"
$ws.Range("C4").Value = "Ror2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.09370099999999999
$ws.Range("H4").Value = 0.281103
$ws.Range("I4").Value = 0.6634967391997054
$ws.Range("J4").Value = 0.6634967391997054
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.159807
$ws.Range("N4").Value = 0.479421
$ws.Range("O4").Value = 0.02005588712602732
$ws.Range("P4").Value = 0.02005588712602732
$ws.Range("Q4").Value = 0.014974075707
$ws.Range("R4").Value = 0.134766681363
$ws.Range("S4").Value = 0.01330701570987648
$ws.Range("T4").Value = 0.01330701570987648

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt1"
$ws.Range("C5").Value = "Ror2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.023296
$ws.Range("H5").Value = 0.069888
$ws.Range("I5").Value = 0.1649589656075852
$ws.Range("J5").Value = 0.1649589656075852
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.11549
$ws.Range("N5").Value = 0.34647
$ws.Range("O5").Value = 0.01449407350231777
$ws.Range("P5").Value = 0.01449407350231777
$ws.Range("Q5").Value = 0.00269045504
$ws.Range("R5").Value = 0.02421409536
$ws.Range("S5").Value = 0.002390927372382648
$ws.Range("T5").Value = 0.002390927372382648

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt1"
$ws.Range("C6").Value = "Ror2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.023296
$ws.Range("H6").Value = 0.069888
$ws.Range("I6").Value = 0.1649589656075852
$ws.Range("J6").Value = 0.1649589656075852
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 7.692787333333334
$ws.Range("N6").Value = 23.078362
$ws.Range("O6").Value = 0.9654500393716549
$ws.Range("P6").Value = 0.965450039371655
$ws.Range("Q6").Value = 0.1792111737173334
$ws.Range("R6").Value = 1.612900563456
$ws.Range("S6").Value = 0.1592596398405506
$ws.Range("T6").Value = 0.1592596398405506

# Row 7: FAPs -> MuSCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt1"
$ws.Range("C7").Value = "Ror2"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.023296
$ws.Range("H7").Value = 0.069888
$ws.Range("I7").Value = 0.1649589656075852
$ws.Range("J7").Value = 0.1649589656075852
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.159807
$ws.Range("N7").Value = 0.479421
$ws.Range("O7").Value = 0.02005588712602732
$ws.Range("P7").Value = 0.02005588712602732
$ws.Range("Q7").Value = 0.003722863872
$ws.Range("R7").Value = 0.033505774848
$ws.Range("S7").Value = 0.003308398394651951
$ws.Range("T7").Value = 0.003308398394651951

# Row 8: MuSCs -> ECs
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Wnt1"
$ws.Range("C8").Value = "Ror2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.024226
$ws.Range("H8").Value = 0.072678
$ws.Range("I8").Value = 0.1715442951927094
$ws.Range("J8").Value = 0.1715442951927094
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.11549
$ws.Range("N8").Value = 0.34647
$ws.Range("O8").Value = 0.01449407350231777
$ws.Range("P8").Value = 0.01449407350231777
$ws.Range("Q8").Value = 0.00279786074
$ws.Range("R8").Value = 0.02518074666
$ws.Range("S8").Value = 0.002486375623426426
$ws.Range("T8").Value = 0.002486375623426427

# Row 9: MuSCs -> FAPs
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Wnt1"
$ws.Range("C9").Value = "Ror2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.024226
$ws.Range("H9").Value = 0.072678
$ws.Range("I9").Value = 0.1715442951927094
$ws.Range("J9").Value = 0.1715442951927094
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 7.692787333333334
$ws.Range("N9").Value = 23.078362
$ws.Range("O9").Value = 0.9654500393716549
$ws.Range("P9").Value = 0.965450039371655
$ws.Range("Q9").Value = 0.1863654659373334
$ws.Range("R9").Value = 1.677289193436
$ws.Range("S9").Value = 0.1656174465477841
$ws.Range("T9").Value = 0.1656174465477841

# Row 10: MuSCs -> MuSCs
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Wnt1"
$ws.Range("C10").Value = "Ror2"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.024226
$ws.Range("H10").Value = 0.072678
$ws.Range("I10").Value = 0.1715442951927094
$ws.Range("J10").Value = 0.1715442951927094
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.159807
$ws.Range("N10").Value = 0.479421
$ws.Range("O10").Value = 0.02005588712602732
$ws.Range("P10").Value = 0.02005588712602732
$ws.Range("Q10").Value = 0.003871484382
$ws.Range("R10").Value = 0.034843359438
$ws.Range("S10").Value = 0.003440473021498892
$ws.Range("T10").Value = 0.003440473021498892

